# Generate Report for Handoff
#
# The localization-status report is regenerated: the "ca7626e9..." file is
# still "In Translation" (handoff data unchanged), while the
# "389a5637..." file has progressed to "Ready for handoff" with a new
# handoff datetime. The regenerated report also re-sorted the two data
# rows (ca7626e9 first, 389a5637 second) on every sheet, but reused the
# previously-created hyperlink relationships (so after the refresh, the
# A2/A3 - etc. hyperlink's underlying target URL corresponds to the row's
# PREVIOUS occupant, not its new one - matching the source report tool's
# behavior of refreshing cell text/hyperlink display without
# re-resolving relationship targets).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Hyperlinks.Delete()

$ws.Range("A2").Value = "ca7626e9-77ca-429b-a63b-133c07e27a8f.md"
$ws.Range("B2").Value = "In Translation"
$ws.Range("C2").Value = "In Translation"
$ws.Range("D2").Value = "2016-16-13 08:16:31"

$ws.Range("A3").Value = "389a5637-3b6f-4151-b8be-ee382053b784.md"
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("D3").Value = "2016-18-13 08:18:54"

$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/f90cbe085f25bcba1152df9a7d92a768ec8b7f01/e2e/389a5637-3b6f-4151-b8be-ee382053b784.md", "", "", "ca7626e9-77ca-429b-a63b-133c07e27a8f.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/f90cbe085f25bcba1152df9a7d92a768ec8b7f01/e2e/ca7626e9-77ca-429b-a63b-133c07e27a8f.md", "", "", "389a5637-3b6f-4151-b8be-ee382053b784.md") | Out-Null

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Hyperlinks.Delete()

$ws.Range("A2").Value = "ca7626e9-77ca-429b-a63b-133c07e27a8f.md"
$ws.Range("B2").Value = ".md"
$ws.Range("C2").Value = "In Translation"
$ws.Range("D2").Value = "ca7626e9-77ca-429b-a63b-133c07e27a8f.40f0fc62222086691ec8629a8907ff29f7d866af.zh-cn.xlf"
$ws.Range("E2").Value = "2016-03-13 08:14:59"
$ws.Range("H2").Value = "0001-01-01 00:00:00"
$ws.Range("I2").Value = "Include"

$ws.Range("A3").Value = "389a5637-3b6f-4151-b8be-ee382053b784.md"
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("D3").Value = "389a5637-3b6f-4151-b8be-ee382053b784.55ce7f4660183f66c14df13c1787ae652a9d2bd4.zh-cn.xlf"
$ws.Range("E3").Value = "2016-03-13 08:18:50"
$ws.Range("H3").Value = "0001-01-01 00:00:00"
$ws.Range("I3").Value = "Include"

$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/f90cbe085f25bcba1152df9a7d92a768ec8b7f01/e2e/389a5637-3b6f-4151-b8be-ee382053b784.md", "", "", "ca7626e9-77ca-429b-a63b-133c07e27a8f.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/f90cbe085f25bcba1152df9a7d92a768ec8b7f01/e2e/389a5637-3b6f-4151-b8be-ee382053b784.md", "", "", ".md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/63c46481e78c1335bbea8d5fd8787b88571f1a54/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/389a5637-3b6f-4151-b8be-ee382053b784.55ce7f4660183f66c14df13c1787ae652a9d2bd4.zh-cn.xlf", "", "", "ca7626e9-77ca-429b-a63b-133c07e27a8f.40f0fc62222086691ec8629a8907ff29f7d866af.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/f90cbe085f25bcba1152df9a7d92a768ec8b7f01/e2e/ca7626e9-77ca-429b-a63b-133c07e27a8f.md", "", "", "389a5637-3b6f-4151-b8be-ee382053b784.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/f90cbe085f25bcba1152df9a7d92a768ec8b7f01/e2e/ca7626e9-77ca-429b-a63b-133c07e27a8f.md", "", "", ".md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/63c46481e78c1335bbea8d5fd8787b88571f1a54/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/ca7626e9-77ca-429b-a63b-133c07e27a8f.40f0fc62222086691ec8629a8907ff29f7d866af.zh-cn.xlf", "", "", "389a5637-3b6f-4151-b8be-ee382053b784.55ce7f4660183f66c14df13c1787ae652a9d2bd4.zh-cn.xlf") | Out-Null

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Hyperlinks.Delete()

$ws.Range("A2").Value = "ca7626e9-77ca-429b-a63b-133c07e27a8f.md"
$ws.Range("B2").Value = ".md"
$ws.Range("C2").Value = "In Translation"
$ws.Range("D2").Value = "ca7626e9-77ca-429b-a63b-133c07e27a8f.40f0fc62222086691ec8629a8907ff29f7d866af.de-de.xlf"
$ws.Range("E2").Value = "2016-03-13 08:16:31"
$ws.Range("H2").Value = "0001-01-01 00:00:00"
$ws.Range("I2").Value = "Include"

$ws.Range("A3").Value = "389a5637-3b6f-4151-b8be-ee382053b784.md"
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("D3").Value = "389a5637-3b6f-4151-b8be-ee382053b784.55ce7f4660183f66c14df13c1787ae652a9d2bd4.de-de.xlf"
$ws.Range("E3").Value = "2016-03-13 08:18:54"
$ws.Range("H3").Value = "0001-01-01 00:00:00"
$ws.Range("I3").Value = "Include"

$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/f90cbe085f25bcba1152df9a7d92a768ec8b7f01/e2e/389a5637-3b6f-4151-b8be-ee382053b784.md", "", "", "ca7626e9-77ca-429b-a63b-133c07e27a8f.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/f90cbe085f25bcba1152df9a7d92a768ec8b7f01/e2e/389a5637-3b6f-4151-b8be-ee382053b784.md", "", "", ".md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d3b00e3660c370373f6664b980a23361e2ff4d98/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/389a5637-3b6f-4151-b8be-ee382053b784.55ce7f4660183f66c14df13c1787ae652a9d2bd4.de-de.xlf", "", "", "ca7626e9-77ca-429b-a63b-133c07e27a8f.40f0fc62222086691ec8629a8907ff29f7d866af.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/f90cbe085f25bcba1152df9a7d92a768ec8b7f01/e2e/ca7626e9-77ca-429b-a63b-133c07e27a8f.md", "", "", "389a5637-3b6f-4151-b8be-ee382053b784.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/f90cbe085f25bcba1152df9a7d92a768ec8b7f01/e2e/ca7626e9-77ca-429b-a63b-133c07e27a8f.md", "", "", ".md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d3b00e3660c370373f6664b980a23361e2ff4d98/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/ca7626e9-77ca-429b-a63b-133c07e27a8f.40f0fc62222086691ec8629a8907ff29f7d866af.de-de.xlf", "", "", "389a5637-3b6f-4151-b8be-ee382053b784.55ce7f4660183f66c14df13c1787ae652a9d2bd4.de-de.xlf") | Out-Null
